# 自动更新Excel文件 - 模拟经过一天后的库存/剩余天数刷新
# For every data row:
#   - D (总天, total days), E (剩余, remaining days), F (开始时间, start date as yyyyMMdd)
#   - If F is not a parseable yyyyMMdd date, leave the row untouched (can't compute).
#   - If E has already reached 1 (about to run out), "refill": reset E back to the
#     total D, and roll the start date F forward by 7 days.
#   - Otherwise, one day has passed, so decrement E by 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {

    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value()
    $eVal = $eCell.Value()
    $fVal = $fCell.Value()

    if ($null -eq $dVal -or $null -eq $eVal -or $null -eq $fVal) {
        continue
    }

    $fText = [string]([long]$fVal)

    $parsedOk = $true
    try {
        $startDate = [datetime]::ParseExact($fText, "yyyyMMdd", $null)
    } catch {
        $parsedOk = $false
    }

    if (-not $parsedOk) {
        continue
    }

    $remaining = [int]$eVal
    $total = [int]$dVal

    if ($remaining -eq 1) {
        # Item expired - refill: reset remaining days and push start date out a week.
        $newDate = $startDate.AddDays(7)
        $eCell.Value = $total
        $fCell.Value = [long]$newDate.ToString("yyyyMMdd")
    } else {
        # One more day has elapsed.
        $eCell.Value = $remaining - 1
    }
}
